$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: shift column labels - C1/D1/E1 become prediction/rejection-f/max
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C2 becomes the species string (was numeric max), D2 stays the species
# string, E2 becomes a numeric rejection-f score (was the species string)
$ws.Range("C2").Value = "s__CAG-288 sp000437395"
$ws.Range("D2").Value = "s__CAG-288 sp000437395"
$ws.Range("E2").Value = 0.99999999999864

# Row 3
$ws.Range("C3").Value = "s__CAG-288 sp000437395"
$ws.Range("D3").Value = "s__CAG-288 sp000437395"
$ws.Range("E3").Value = 0.9999999999986302

# Row 4
$ws.Range("C4").Value = "s__CAG-288 sp000437395"
$ws.Range("D4").Value = "s__CAG-288 sp000437395"
$ws.Range("E4").Value = 0.9999999999985483

# Row 5
$ws.Range("C5").Value = "s__CAG-288 sp000437395"
$ws.Range("D5").Value = "s__CAG-288 sp000437395"
$ws.Range("E5").Value = 0.999999999998779
